# Add a new "5-jul" data column (W) to the dataframe, mirroring the
# structure of the existing "4-jul" column (V).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column
$ws.Range("W1").Value = "5-jul"

# New column values, row by row (rows 2-18)
$ws.Range("W2").Value  = 0
$ws.Range("W3").Value  = 14.479368061886348
$ws.Range("W4").Value  = 12.815287227932272
$ws.Range("W5").Value  = 27.873531228416816
$ws.Range("W6").Value  = 0
$ws.Range("W7").Value  = 33.880299523869624
$ws.Range("W8").Value  = 16.545371833125643
$ws.Range("W9").Value  = 21.356382587299496
$ws.Range("W10").Value = 26.604011613386586
$ws.Range("W11").Value = 13.460144203373545
$ws.Range("W12").Value = 0
$ws.Range("W13").Value = 13.936386260932005
$ws.Range("W14").Value = 0
$ws.Range("W15").Value = 0
$ws.Range("W16").Value = 14.838542045497856
$ws.Range("W17").Value = 0
$ws.Range("W18").Value = 0

# Reflect the updated view/selection state seen in the authored workbook
$ws.Activate()
$ws.Range("U4").Select()
